# "Put my full name" — split the title run "Saswat Levin" so that the
# middle name "Kuthully" is inserted between "Saswat" and "Levin",
# producing three separate (but identically-formatted) runs as in the
# target diff.

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$titleRange = $p.Range

# Locate the boundary right after "Saswat " (7 characters in) where the
# new middle name needs to be inserted.
$insertAt = $titleRange.Start + 7
$insertionPoint = $d.Range($insertAt, $insertAt)
$insertionPoint.InsertBefore("Kuthully ")

# The freshly inserted text shares identical run formatting with its
# neighbours, so the engine would normally coalesce it back into a single
# run. Toggling a character property off and on again forces the text to
# remain in its own run while leaving the final formatting unchanged.
$newWordRange = $d.Range($insertAt, $insertAt + 9)
$newWordRange.Font.Bold = $false
$newWordRange.Font.Bold = $true
